$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new weekly records were added at the top of the data block (rows 99-100),
# pushing the existing rows 99:180 down to 101:182.
$ws.Rows("99:100").Insert()

# New row 99: Melón, Calameño, Primera — fecha 44546
$ws.Range("A99").Value = 5
$ws.Range("B99").Value = "Macroferia Regional de Talca"
$ws.Range("C99").Value = "Maule"
$ws.Range("D99").Value = 44546
$ws.Range("E99").Value = 7
$ws.Range("F99").Value = 100112027
$ws.Range("G99").Value = "Melón"
$ws.Range("H99").Value = "Calameño"
$ws.Range("I99").Value = "Primera"
$ws.Range("J99").Value = 3000
$ws.Range("K99").Value = 1000
$ws.Range("L99").Value = 1000
$ws.Range("M99").Value = 1000
$ws.Range("N99").Value = "`$/unidad"
$ws.Range("O99").Value = "Región del Maule"
$ws.Range("P99").Value = 1000
$ws.Range("Q99").Value = 1
$ws.Range("R99").Value = "Hortaliza"

# New row 100: Melón, Calameño, Segunda — fecha 44546
$ws.Range("A100").Value = 5
$ws.Range("B100").Value = "Macroferia Regional de Talca"
$ws.Range("C100").Value = "Maule"
$ws.Range("D100").Value = 44546
$ws.Range("E100").Value = 7
$ws.Range("F100").Value = 100112027
$ws.Range("G100").Value = "Melón"
$ws.Range("H100").Value = "Calameño"
$ws.Range("I100").Value = "Segunda"
$ws.Range("J100").Value = 4000
$ws.Range("K100").Value = 800
$ws.Range("L100").Value = 800
$ws.Range("M100").Value = 800
$ws.Range("N100").Value = "`$/unidad"
$ws.Range("O100").Value = "Región del Maule"
$ws.Range("P100").Value = 800
$ws.Range("Q100").Value = 1
$ws.Range("R100").Value = "Hortaliza"
